$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update a few existing daily totals (August 2025 block)
$ws.Range("B3").Value = 30349.56
$ws.Range("B10").Value = 26019
$ws.Range("B11").Value = 14123.9

# Insert a new daily record for August 15, 2025 right after August 14 (row 11),
# pushing the July/June/May blocks down by one row.
$ws.Rows(12).Insert()
$ws.Range("A12").Value = 15
$ws.Range("B12").Value = 49557.28
$ws.Range("C12").Value = 8
$ws.Range("D12").Value = 2025
$ws.Range("E12").Value = "08/2025"
